$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PARAMETERS")

$ws.Range("B2").Value = "VVI"
$ws.Range("C2").Value = "40"
$ws.Range("D2").Value = "120"
$ws.Range("F2").Value = "3.5"
$ws.Range("H2").Value = "0.4"
$ws.Range("J2").Value = "320"
$ws.Range("K2").Value = "250"
$ws.Range("L2").Value = "Med"
$ws.Range("M2").Value = "30"
$ws.Range("N2").Value = "8"
